$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1334.1428
$ws.Range("I129").Value = 1023.8333
$ws.Range("K129").Value = 3071.4999
$ws.Range("M129").Value = 1928.5001
$ws.Range("H132").Value = 8550.666999999999
$ws.Range("I132").Value = 10596.429
$ws.Range("K132").Value = 31789.287
$ws.Range("M132").Value = -29259.287
$ws.Range("H141").Value = 4163.3335
$ws.Range("I141").Value = 3497.5
$ws.Range("J141").Value = 5495
$ws.Range("K141").Value = 10492.5
$ws.Range("L141").Value = 16485
$ws.Range("M141").Value = -5312.5
$ws.Range("N141").Value = -26845

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
$ws.Range("H4").Value = 643.4
$ws.Range("I4").Value = 332
$ws.Range("K4").Value = 332
$ws.Range("M4").Value = -216
$ws.Range("H6").Value = 22222222
$ws.Range("I6").Value = 24000000
$ws.Range("K6").Value = 24000000
$ws.Range("M6").Value = -23999827
$ws.Range("H32").Value = 743.14813
$ws.Range("I32").Value = 743.14813
$ws.Range("K32").Value = 743.14813
$ws.Range("M32").Value = -456.14813
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H74").Value = 1548.5
$ws.Range("I74").Value = 1097
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1097
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -223
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1548.5
$ws.Range("I77").Value = 1097
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 5485
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -1117
$ws.Range("N77").Value = -18736
$ws.Range("H97").Value = 694.0909
$ws.Range("I97").Value = 958
$ws.Range("J97").Value = 474.16666
$ws.Range("K97").Value = 958
$ws.Range("L97").Value = 474.16666
$ws.Range("M97").Value = -462
$ws.Range("N97").Value = -1466.16666
$ws.Range("H110").Value = 457.625
$ws.Range("I110").Value = 457.625
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 457.625
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1587.375
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 1776.2222
$ws.Range("I132").Value = 1776.2222
$ws.Range("K132").Value = 5328.6666
$ws.Range("M132").Value = -2798.6666

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 4414.143
$ws.Range("I22").Value = 4733.1665
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 4733.1665
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -4560.1665
$ws.Range("N22").Value = -2846
$ws.Range("H26").Value = 14808.625
$ws.Range("I26").Value = 15638.429
$ws.Range("K26").Value = 15638.429
$ws.Range("M26").Value = -15346.429
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H96").Value = 13570.4
$ws.Range("I96").Value = 13570.4
$ws.Range("K96").Value = 13570.4
$ws.Range("M96").Value = -10824.4
$ws.Range("H105").Value = 1211.8823
$ws.Range("I105").Value = 1109.9286
$ws.Range("K105").Value = 1109.9286
$ws.Range("M105").Value = 637.0714
$ws.Range("H134").Value = 2238.5
$ws.Range("I134").Value = 1798.125
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5394.375
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -2859.375
$ws.Range("N134").Value = -17070

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2191.25
$ws.Range("I22").Value = 1959.6
$ws.Range("J22").Value = 2356.7144
$ws.Range("K22").Value = 1959.6
$ws.Range("L22").Value = 2356.7144
$ws.Range("M22").Value = -1609.6
$ws.Range("N22").Value = -3056.7144
$ws.Range("H31").Value = 5047.9756
$ws.Range("I31").Value = 2877.8215
$ws.Range("K31").Value = 2877.8215
$ws.Range("M31").Value = -2582.8215
$ws.Range("H34").Value = 5047.9756
$ws.Range("I34").Value = 2877.8215
$ws.Range("K34").Value = 2877.8215
$ws.Range("M34").Value = -2675.8215
$ws.Range("H62").Value = 8650.666999999999
$ws.Range("I62").Value = 8200
$ws.Range("K62").Value = 8200
$ws.Range("M62").Value = -7576
$ws.Range("H65").Value = 8650.666999999999
$ws.Range("I65").Value = 8200
$ws.Range("K65").Value = 41000
$ws.Range("M65").Value = -37880
$ws.Range("H93").Value = 3215.6667
$ws.Range("I93").Value = 3215.6667
$ws.Range("K93").Value = 3215.6667
$ws.Range("M93").Value = -1343.6667
$ws.Range("H103").Value = 13391.6
$ws.Range("I103").Value = 13391.6
$ws.Range("K103").Value = 13391.6
$ws.Range("M103").Value = -12219.6

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 562.7692
$ws.Range("I2").Value = 393.16666
$ws.Range("J2").Value = 708.1429000000001
$ws.Range("K2").Value = 393.16666
$ws.Range("L2").Value = 708.1429000000001
$ws.Range("M2").Value = -280.16666
$ws.Range("N2").Value = -934.1429000000001
$ws.Range("H23").Value = 404
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 600
$ws.Range("K23").Value = 12
$ws.Range("L23").Value = 600
$ws.Range("M23").Value = 211
$ws.Range("N23").Value = -1046
$ws.Range("H24").Value = 46295.832
$ws.Range("J24").Value = 46295.832
$ws.Range("L24").Value = 46295.832
$ws.Range("N24").Value = -46641.832
$ws.Range("H43").Value = 17349.834
$ws.Range("J43").Value = 17349.834
$ws.Range("L43").Value = 17349.834
$ws.Range("N43").Value = -17651.834
$ws.Range("H46").Value = 11250
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2344
$ws.Range("H57").Value = 15333.333
$ws.Range("I57").Value = 5500
$ws.Range("K57").Value = 5500
$ws.Range("M57").Value = -4680
$ws.Range("H86").Value = 16304
$ws.Range("J86").Value = 16304
$ws.Range("L86").Value = 16304
$ws.Range("N86").Value = -18676
$ws.Range("H89").Value = 16304
$ws.Range("J89").Value = 16304
$ws.Range("L89").Value = 48912
$ws.Range("N89").Value = -60768
$ws.Range("H97").Value = 1400
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 2056.647
$ws.Range("I102").Value = 1427.3077
$ws.Range("K102").Value = 1427.3077
$ws.Range("M102").Value = 194.6922999999999
$ws.Range("H126").Value = 1799.8
$ws.Range("I126").Value = 1749.75
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5249.25
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2779.25
$ws.Range("N126").Value = -10940
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5931.125
$ws.Range("I46").Value = 483
$ws.Range("J46").Value = 9200
$ws.Range("K46").Value = 483
$ws.Range("L46").Value = 9200
$ws.Range("M46").Value = -295
$ws.Range("N46").Value = -9576
$ws.Range("H106").Value = 11333.8
$ws.Range("J106").Value = 11333.8
$ws.Range("L106").Value = 11333.8
$ws.Range("N106").Value = -13857.8
$ws.Range("H122").Value = 2560.6
$ws.Range("I122").Value = 2560.6
$ws.Range("K122").Value = 7681.799999999999
$ws.Range("M122").Value = -5231.799999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1055.75
$ws.Range("I81").Value = 1055.75
$ws.Range("K81").Value = 2111.5
$ws.Range("M81").Value = -1050.5
$ws.Range("H84").Value = 1055.75
$ws.Range("I84").Value = 1055.75
$ws.Range("K84").Value = 10557.5
$ws.Range("M84").Value = -5253.5
$ws.Range("H105").Value = 12666.667
$ws.Range("J105").Value = 12666.667
$ws.Range("L105").Value = 12666.667
$ws.Range("N105").Value = -19654.667
$ws.Range("H107").Value = 395.46667
$ws.Range("I107").Value = 352.2857
$ws.Range("K107").Value = 1056.8571
$ws.Range("M107").Value = 863.1428999999998
$ws.Range("H122").Value = 3248
$ws.Range("I122").Value = 2097.6
$ws.Range("K122").Value = 6292.799999999999
$ws.Range("M122").Value = -3842.799999999999
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 1000
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
